# Backup codes rotation: the first three codes (A2:A4) have been consumed,
# so pull the next three available codes (A8:A10) up into their place and
# clear out the now-used-up rows 8:10 (rows 11:12 stay put as the
# remaining reserve). Finally move the active selection to A5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = $ws.Range("A8").Value()
$ws.Range("A3").Value = $ws.Range("A9").Value()
$ws.Range("A4").Value = $ws.Range("A10").Value()

$ws.Range("A8:A10").ClearContents()

$ws.Range("A5").Select()
